# Manage Claims Email Update
# Two new claims (57794958 / 57794961) are recorded, replacing the two
# oldest entries shown on the "Input", "ClaimDetail" and "ShipmentInformation"
# sheets (rows 2 & 3 act as a rolling "most recent claims" log).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# ShipmentInformation sheet
# ---------------------------------------------------------------------
$wsShip = $wb.Worksheets.Item("ShipmentInformation")

$wsShip.Range("C2").NumberFormat = "@"
$wsShip.Range("C2").Value = "PickUp90"

$wsShip.Range("K2").NumberFormat = "@"
$wsShip.Range("K2").Value = "DropOff660"

# ---------------------------------------------------------------------
# Input sheet
# ---------------------------------------------------------------------
$wsInput = $wb.Worksheets.Item("Input")

# Row 2 -> newest claim (57794958)
$wsInput.Range("B2").NumberFormat = "@"
$wsInput.Range("B2").Value = "10-15-2021"

$wsInput.Range("T2").NumberFormat = "@"
$wsInput.Range("T2").Value = "57794958"

$wsInput.Range("U2").NumberFormat = "@"
$wsInput.Range("U2").Value = "`$688.73"

$wsInput.Range("W2").NumberFormat = "@"
$wsInput.Range("W2").Value = "FCT898607587350544384"

$wsInput.Range("X2").NumberFormat = "@"
$wsInput.Range("X2").Value = "FCTEST1003764"

$wsInput.Range("Y2").NumberFormat = "@"
$wsInput.Range("Y2").Value = "`$27.27"

# Row 3 -> next claim (57794961)
$wsInput.Range("B3").NumberFormat = "@"
$wsInput.Range("B3").Value = "10-15-2021"

$wsInput.Range("T3").NumberFormat = "@"
$wsInput.Range("T3").Value = "57794961"

$wsInput.Range("W3").NumberFormat = "@"
$wsInput.Range("W3").Value = "FCT898617792033456128"

$wsInput.Range("X3").NumberFormat = "@"
$wsInput.Range("X3").Value = "FCTEST1003765"

# ---------------------------------------------------------------------
# ClaimDetail sheet
# ---------------------------------------------------------------------
$wsClaim = $wb.Worksheets.Item("ClaimDetail")

# Row 2 -> newest claim (57794958), status moved from Initiated to Filed
$wsClaim.Range("A2").NumberFormat = "@"
$wsClaim.Range("A2").Value = "57794958"

$wsClaim.Range("B2").NumberFormat = "@"
$wsClaim.Range("B2").Value = "10-15-2021"

$wsClaim.Range("C2").NumberFormat = "@"
$wsClaim.Range("C2").Value = "Filed"

# Row 3 -> next claim (57794961), stays Filed
$wsClaim.Range("A3").NumberFormat = "@"
$wsClaim.Range("A3").Value = "57794961"

$wsClaim.Range("B3").NumberFormat = "@"
$wsClaim.Range("B3").Value = "10-15-2021"

$wsClaim.Range("C3").NumberFormat = "@"
$wsClaim.Range("C3").Value = "Filed"
